# Apply the data refresh edit described by the commit:
# "Upload excel files with prices"
#
# 1) All row timestamps (column O, rows 2-397) are bumped from the previous
#    crawl time to the new crawl time.
# 2) A handful of rows had their ratingAmount (column D) - and in a few
#    cases ratingValue (column E) - updated to reflect newly crawled values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2023-01-13 06:49:14"
$newTimestamp = "2023-01-13 12:56:17"

$lastRow = 397

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Column D (ratingAmount) updates
$dUpdates = @{
    5   = 25
    8   = 27
    34  = 13
    39  = 15
    44  = 38
    48  = 9
    52  = 16
    73  = 28
    76  = 18
    78  = 18
    89  = 17
    92  = 10
    95  = 64
    100 = 7
    122 = 7
    130 = 6
    133 = 15
    156 = 13
    168 = 19
    187 = 20
    225 = 24
    280 = 6
    384 = 41
}

foreach ($row in $dUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
}

# Column E (ratingValue) updates
$eUpdates = @{
    78  = 4.5
    92  = 3.5
    122 = 4.5
    187 = 5
    280 = 3.5
}

foreach ($row in $eUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $eUpdates[$row]
}
